$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values that changed (use a text formula then harden it to a
# static value via copy/paste-special so the cell keeps its original style
# and is stored as plain text rather than being auto-coerced to a number
# or date).
$ws.Range("A2").Formula = '="201297"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("C2").Formula = '="11/10/2025"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)

# Remove the now-obsolete rows 3-5, shifting cells up (matches the diff's
# row deletion of the three extra excuse entries).
$ws.Range("A3:F5").Delete()
